$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "60.194.43"
$ws.Cells.Item(2, 5).Value = "  -4.22%  "
$ws.Cells.Item(3, 4).Value = "2.903.30"
$ws.Cells.Item(3, 5).Value = "  -3.51%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "527.71"
$ws.Cells.Item(5, 5).Value = "  -5.16%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "141.60"
$ws.Cells.Item(6, 5).Value = "  -7.65%  "
$ws.Cells.Item(7, 5).Value = "  +0.06%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.549"
$ws.Cells.Item(8, 5).Value = "  -2.28%  "
$ws.Cells.Item(9, 4).Value = "2.903.65"
$ws.Cells.Item(9, 5).Value = "  -3.68%  "
$ws.Cells.Item(10, 5).Value = "  -5.03%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.85"
$ws.Cells.Item(11, 5).Value = "  -8.57%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.353"
$ws.Cells.Item(12, 5).Value = "  -3.20%  "
$ws.Cells.Item(13, 4).Value = "3.409.64"
$ws.Cells.Item(13, 5).Value = "  -3.55%  "
$ws.Cells.Item(14, 5).Value = "  +1.27%  "
$ws.Cells.Item(15, 4).Value = "60.423.25"
$ws.Cells.Item(15, 5).Value = "  -4.01%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "22.60"
$ws.Cells.Item(16, 5).Value = "  -5.48%  "
$ws.Cells.Item(17, 4).Value = "2.902.18"
$ws.Cells.Item(17, 5).Value = "  -3.65%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.0000139"
$ws.Cells.Item(18, 5).Value = "  -6.73%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "4.90"
$ws.Cells.Item(19, 5).Value = "  -3.93%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.47"
$ws.Cells.Item(20, 5).Value = "  -3.43%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "359.37"
$ws.Cells.Item(21, 5).Value = "  -9.19%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.57"
$ws.Cells.Item(22, 5).Value = "  -0.72%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "1.00"
$ws.Cells.Item(23, 5).Value = "  -0.01%  "
$ws.Cells.Item(24, 5).Value = "  -1.93%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "63.32"
$ws.Cells.Item(25, 5).Value = "  -2.72%  "
$ws.Cells.Item(26, 4).Value = "3.027.28"
$ws.Cells.Item(26, 5).Value = "  -3.75%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.446"
$ws.Cells.Item(27, 5).Value = "  -4.12%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.174"
$ws.Cells.Item(28, 5).Value = "  -7.21%  "
$ws.Cells.Item(29, 5).Value = "  +0.55%  "
$ws.Cells.Item(30, 5).Value = "  -11.35%  "
$ws.Cells.Item(31, 4).Value = "0.0₃0846"
$ws.Cells.Item(31, 5).Value = "  -12.66%  "
$ws.Cells.Item(32, 5).Value = "  -0.02%  "
$ws.Cells.Item(33, 5).Value = "  -5.59%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "19.47"
$ws.Cells.Item(34, 5).Value = "  -4.84%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "149.71"
$ws.Cells.Item(35, 5).Value = "  -6.21%  "
$ws.Cells.Item(36, 5).Value = "  -8.82%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.50"
$ws.Cells.Item(37, 5).Value = "  -8.86%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.988"
$ws.Cells.Item(38, 5).Value = "  -9.74%  "
$ws.Cells.Item(39, 5).Value = "  -8.46%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "37.55"
$ws.Cells.Item(40, 5).Value = "  -0.14%  "
$ws.Cells.Item(41, 4).Value = "2.330.47"
$ws.Cells.Item(41, 5).Value = "  -6.67%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.44"
$ws.Cells.Item(42, 5).Value = "  -8.47%  "
$ws.Cells.Item(43, 5).Value = "  -6.77%  "
$ws.Cells.Item(44, 5).Value = "  -3.54%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "20.69"
$ws.Cells.Item(45, 5).Value = "  -8.18%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0566"
$ws.Cells.Item(46, 5).Value = "  -4.82%  "
$ws.Cells.Item(47, 5).Value = "  -0.01%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "4.79"
$ws.Cells.Item(48, 5).Value = "  -4.53%  "
$ws.Cells.Item(49, 5).Value = "  -1.22%  "
$ws.Cells.Item(50, 5).Value = "  -6.11%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0922"
$ws.Cells.Item(51, 5).Value = "  -2.33%  "
